# Updated capital structure database
# Applies the refreshed metrics for the two Kenya "Investments & Asset
# Management" rows (rows 2 and 3) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2, 3) {
    # historical_growth_revenue_last_5_years replaces the old D/E pair;
    # historical_growth_net_income_last_5_years (E) is dropped entirely.
    $ws.Range("D$r").Value = -0.0788
    $ws.Range("E$r").ClearContents()

    $ws.Range("G$r").Value = 1.214477211796247
    $ws.Range("H$r").Value = 1.214477211796247
    $ws.Range("I$r").Value = 0.06899016979445934
    $ws.Range("J$r").Value = 0.06899016979445934
    $ws.Range("K$r").Value = -18
    $ws.Range("L$r").Value = -0.160857908847185
    $ws.Range("M$r").Value = 0
    $ws.Range("N$r").Value = 0
    $ws.Range("O$r").Value = 0
    $ws.Range("P$r").Value = 0
    $ws.Range("Q$r").Value = 0
    $ws.Range("R$r").Value = 0
    # S stays 0 (unchanged); buybacks_cash_returned (T) is dropped entirely.
    $ws.Range("T$r").ClearContents()

    $ws.Range("U$r").Value = 81.5
    $ws.Range("V$r").Value = 0.8333333333333334
    $ws.Range("W$r").Value = -0.0377992440151197
    $ws.Range("X$r").Value = 0.164420107173175
    $ws.Range("Y$r").Value = -0.2022193511882947
    $ws.Range("Z$r").Value = 0.2185546875
    $ws.Range("AA$r").Value = 0.015078125
    $ws.Range("AB$r").Value = 0.08012104447491444
    $ws.Range("AC$r").Value = -0.06504291947491445
    $ws.Range("AD$r").Value = 252.6
    # AE stays 0 (unchanged)
    $ws.Range("AF$r").Value = 252.6
    $ws.Range("AG$r").Value = 171.1
    $ws.Range("AH$r").Value = 0.7208904109589042
    $ws.Range("AI$r").Value = 0.3543771043771044
    $ws.Range("AJ$r").Value = 0.6362960208255858
    $ws.Range("AK$r").Value = 0.2710280373831776
    $ws.Range("AL$r").Value = 16.7
    $ws.Range("AM$r").Value = 16.7
    $ws.Range("AN$r").Value = 20.53658536585366
    $ws.Range("AO$r").Value = 0.4622754491017964
    $ws.Range("AP$r").Value = 13.91056910569105
    $ws.Range("AQ$r").Value = 0.4622754491017964
}
